$d = $word.ActiveDocument

# Locate the target paragraph ("The application has a default account ...")
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*default account (admin account)*") {
        $targetPara = $p
        break
    }
}

$pStart = $targetPara.Range.Start
$pEnd = $targetPara.Range.End

# Find the email run ("admin@FUNewsManagementSystem.org") to bookmark as OLE_LINK1
$emailRange = $d.Range($pStart, $pEnd)
$emailRange.Find.Execute("admin@FUNewsManagementSystem.org", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("OLE_LINK1", $emailRange) | Out-Null

# Find the password value together with its surrounding quotes ("@@abc123@@") to
# bookmark as OLE_LINK2, extended one character on each side to capture the
# opening/closing curly quotes and the trailing space after the closing quote.
$pwRange = $d.Range($pStart, $pEnd)
$pwRange.Find.Execute("@@abc123@@", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmStart = $pwRange.Start - 1
$bmEnd = $pwRange.End + 2
$pwBookmarkRange = $d.Range($bmStart, $bmEnd)
$d.Bookmarks.Add("OLE_LINK2", $pwBookmarkRange) | Out-Null
